$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 2

# Remove row 3 entirely (data + dimension shrinks)
$ws.Range("A3:B3").Delete()
